$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated refrigerant water temperatures (refr in / refr out)
$ws.Range("B2").Value = -62
$ws.Range("B3").Value = -35

# Updated Ucond / Ureboil base values
$ws.Range("B6").Value = 80
$ws.Range("B7").Value = 80

# Update the active selection to match the target state
$ws.Range("E12").Select()
